# Daily cryptos price/volume refresh (GitHub Actions scheduled update).
# The source diff only touches column D (Price) and column E (Volume(1h))
# on Sheet1, row by row; everything else (coin name, link, row order) is
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.621.16"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "'1.870.27"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'324.63"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "'0.4627"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "'0.3879"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.07872"
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'0.9752"
$ws.Range("E10").Value = "  +1.81%  "
$ws.Range("D11").Value = "'21.96"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "'1.887.48"
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "'6.994"
$ws.Range("E13").Value = "  +1.49%  "
$ws.Range("D14").Value = "'5.700"
$ws.Range("E14").Value = "  +0.63%  "
$ws.Range("D15").Value = "'0.06950"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "'88.14"
$ws.Range("E16").Value = "  +1.14%  "
$ws.Range("D17").Value = "'1.006"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("E18").Value = "  +1.15%  "
$ws.Range("D19").Value = "'16.82"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").Value = "'1.004"
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'28.620.46"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("D22").Value = "'5.275"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("D23").Value = "'11.01"
$ws.Range("E23").Value = "  +0.48%  "
$ws.Range("D24").Value = "'2.115"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("D25").Value = "'2.111.82"
$ws.Range("E25").Value = "  +3.16%  "
$ws.Range("D26").Value = "'152.58"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'19.24"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("D29").Value = "'1.985"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "'119.32"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'0.09338"
$ws.Range("E31").Value = "  +0.91%  "
$ws.Range("D32").Value = "'0.9164"
$ws.Range("E32").Value = "  -1.72%  "
$ws.Range("D33").Value = "'5.263"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "'1.337"
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("D35").Value = "'3.330"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").Value = "'0.05789"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "'0.02104"
$ws.Range("E37").Value = "  -1.68%  "
$ws.Range("D38").Value = "'1.154"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'7.736"
$ws.Range("E39").Value = "  -0.26%  "
$ws.Range("D40").Value = "'0.5625"
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "'0.1785"
$ws.Range("E41").Value = "  +1.41%  "
$ws.Range("D42").Value = "'9.779"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "'0.07181"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").Value = "'11.74"
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("D45").Value = "'0.5304"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("D46").Value = "'2.162"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("D47").Value = "'1.139"
$ws.Range("E47").Value = "  +2.58%  "
$ws.Range("D48").Value = "'1.829"
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("D49").Value = "'112.84"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'2.408"
$ws.Range("E50").Value = "  +3.63%  "
$ws.Range("D51").Value = "'1.004"
$ws.Range("E51").Value = "  +0.34%  "
